$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Fill in "Hours Burnt" (column F) for rows 13-17 so that the
# "Remaining Hours" formula in column G (=E-F) recalculates to 0.
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 2

# Update the view state: scroll so column C is left-most visible and
# select E18 as the active cell.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E18").Select()
